$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) "Test Writer:" line - credit a second author, with spell-check marks
#    around the two new (unrecognised) proper nouns, joined with "& " to
#    the original "Sarmad Butti".
# -----------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("Sarmad Butti", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find 'Sarmad Butti' (Test Writer line)"
}

$xmlTestWriter = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Test Writer:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Edgard</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Musafiri</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> &amp; </w:t></w:r><w:r><w:t>Sarmad Butti</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
[void]$r1.InsertXML($xmlTestWriter)

# -----------------------------------------------------------------------
# 2) "Setup:" line - split the leading word "Attaching" -> "Connecting"
#    into its own run, with a (re-dropped) _GoBack bookmark right after
#    it, same place Word leaves the cursor after an edit.
# -----------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("Attaching a battery to the device and start testing", `
                            $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the 'Attaching a battery...' setup sentence"
}

$xmlSetup = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Connecting</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> a battery to the device and start testing </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
[void]$r2.InsertXML($xmlSetup)

# -----------------------------------------------------------------------
# 3) The stray trailing paragraph (after the table) used to carry the
#    _GoBack bookmark; now that the bookmark lives next to "Connecting"
#    it must be removed from here, leaving a plain empty paragraph.
# -----------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

Write-Host "Applied: co-author credit, Setup wording tweak, _GoBack bookmark relocated."
